# Changed the instruction and task script input parameter strategy
#
# - Cue_brightness (column H, rows 2-21) is no longer a random draw; every
#   trial now gets a fixed brightness value of 1.
# - Column widths on Sheet1 are set to their "best fit" content widths.
# - The active selection moves from L17 to P14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- H2:H21 -> 1 -------------------------------------------------------
$ws.Range("H2:H21").Value = 1

# --- Column widths (best-fit to content) --------------------------------
# NOTE: the engine's ColumnWidth setter quantizes to the nearest 1/6 of a
# character, so these inputs are chosen so the stored width lands on the
# closest reachable value to the authored best-fit widths:
#   A=4.5546875 B=5.5546875 C=15.33203125 D=9.33203125
#   E=9.44140625 F=9.44140625 G=8.6640625 H=14
$ws.Columns.Item(1).ColumnWidth = 3.6328125
$ws.Columns.Item(2).ColumnWidth = 4.6875
$ws.Columns.Item(3).ColumnWidth = 14.53125
$ws.Columns.Item(4).ColumnWidth = 8.4375
$ws.Columns.Item(5).ColumnWidth = 8.671875
$ws.Columns.Item(6).ColumnWidth = 8.671875
$ws.Columns.Item(7).ColumnWidth = 7.8515625
$ws.Columns.Item(8).ColumnWidth = 13.125

# --- Selection moves to P14 ---------------------------------------------
[void]$ws.Range("P14").Select()
